# "added range argument to read_excel" - adds a new "position" worksheet
# (the same a/b/c0/c1/c2 table as the "3d" sheet, but offset to D3:H9 so the
# new range argument of read_excel has something non-A1-anchored to read)
# plus a couple of leftover selection tweaks on existing sheets.

$wb = $excel.ActiveWorkbook

# --- 3d: selection becomes the whole used range (no explicit active cell) ---
$ws3d = $wb.Worksheets.Item("3d")
$ws3d.Range("A1:E7").Select()

# --- int_labels: selection moves to H15 ---
$intLabels = $wb.Worksheets.Item("int_labels")
$intLabels.Range("H15").Select()

# --- add the new "position" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$position = $wb.Worksheets.Add($null, $lastSheet)
$position.Name = "position"

# header row (row 3, columns D:H)
$position.Cells.Item(3, 4).Value = "a"
$position.Cells.Item(3, 5).Value = "b\c"
$position.Cells.Item(3, 6).Value = "c0"
$position.Cells.Item(3, 7).Value = "c1"
$position.Cells.Item(3, 8).Value = "c2"

# data rows 4-9, columns D:H
$position.Cells.Item(4, 4).Value = 1
$position.Cells.Item(4, 5).Value = "b0"
$position.Cells.Item(4, 6).Value = 0
$position.Cells.Item(4, 7).Value = 1
$position.Cells.Item(4, 8).Value = 2

$position.Cells.Item(5, 4).Value = 1
$position.Cells.Item(5, 5).Value = "b1"
$position.Cells.Item(5, 6).Value = 3
$position.Cells.Item(5, 7).Value = 4
$position.Cells.Item(5, 8).Value = 5

$position.Cells.Item(6, 4).Value = 2
$position.Cells.Item(6, 5).Value = "b0"
$position.Cells.Item(6, 6).Value = 6
$position.Cells.Item(6, 7).Value = 7
$position.Cells.Item(6, 8).Value = 8

$position.Cells.Item(7, 4).Value = 2
$position.Cells.Item(7, 5).Value = "b1"
$position.Cells.Item(7, 6).Value = 9
$position.Cells.Item(7, 7).Value = 10
$position.Cells.Item(7, 8).Value = 11

$position.Cells.Item(8, 4).Value = 3
$position.Cells.Item(8, 5).Value = "b0"
$position.Cells.Item(8, 6).Value = 12
$position.Cells.Item(8, 7).Value = 13
$position.Cells.Item(8, 8).Value = 14

$position.Cells.Item(9, 4).Value = 3
$position.Cells.Item(9, 5).Value = "b1"
$position.Cells.Item(9, 6).Value = 15
$position.Cells.Item(9, 7).Value = 16
$position.Cells.Item(9, 8).Value = 17

# leave the selection on H9 of the new (now active/tab-selected) sheet
$position.Range("H9").Select()

Write-Output "position sheet added"
